$d = $word.ActiveDocument

# --- Part 1: merge "Stage 2 - Semi" + " " runs into "Stage 2 - " and drop the
#     bookmark that currently sits between them. Find/Replace naturally
#     merges the matched runs into a single run with the replacement text.
$null = $d.Content.Find.Execute(
    "Stage 2 – Semi ", $false, $false, $false, $false, $false,
    $true, 1, $false, "Stage 2 – ", 2)

# The old "_GoBack" bookmark (formerly sitting between the two runs we just
# merged) is gone now that its surrounding text was replaced. Make sure no
# stray bookmark with that name remains before we re-add it elsewhere.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Part 2: (re)create the "_GoBack" bookmark at the very end of the last
#     cell of the table (after "Not resign"), matching the target markup
#     where the bookmark lives inside <w:tr> right after the last <w:tc>.
#
# The engine's Bookmarks.Add falls back to the start of the document when
# given a collapsed Range that sits exactly on the trailing edge of a run
# (i.e. Start == the run's End). To avoid that, we temporarily extend the
# cell's text with a placeholder run, add the bookmark at the boundary
# between the real text and the placeholder (which now resolves correctly
# because it's no longer the run's trailing edge), and then delete the
# placeholder text, leaving the bookmark exactly where we want it.
$table = $d.Tables.Item(1)
$lastRow = $table.Rows.Item($table.Rows.Count)
$lastCell = $lastRow.Cells.Item($lastRow.Cells.Count)

$endRange = $lastCell.Range.Duplicate
$endRange.Collapse(0)
$insertPos = $endRange.Start
$endRange.InsertAfter("ZZ_TMP_ZZ")

$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($insertPos, $insertPos + 9)
$placeholder.Delete()
